$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 '26.908.73'
Set-TextCell 2 5 '  +2.22%  '

# Row 3
Set-TextCell 3 4 '1.809.95'
Set-TextCell 3 5 '  +2.92%  '

# Row 4
Set-TextCell 4 5 '  +0.60%  '

# Row 5
Set-TextCell 5 4 '313.17'
Set-TextCell 5 5 '  +3.68%  '

# Row 6
Set-TextCell 6 4 '1.007'
Set-TextCell 6 5 '  +0.54%  '

# Row 7
Set-TextCell 7 4 '0.4288'
Set-TextCell 7 5 '  +0.14%  '

# Row 8
Set-TextCell 8 5 '  +2.06%  '

# Row 9
Set-TextCell 9 4 '0.07257'
Set-TextCell 9 5 '  +3.19%  '

# Row 10
Set-TextCell 10 4 '0.8647'
Set-TextCell 10 5 '  +4.09%  '

# Row 11
Set-TextCell 11 4 '2.048.90'
Set-TextCell 11 5 '  +18.44%  '

# Row 12
Set-TextCell 12 5 '  +5.44%  '

# Row 13
Set-TextCell 13 4 '6.620'
Set-TextCell 13 5 '  +4.30%  '

# Row 14
Set-TextCell 14 4 '5.394'
Set-TextCell 14 5 '  +3.56%  '

# Row 15
Set-TextCell 15 4 '0.06929'
Set-TextCell 15 5 '  +2.15%  '

# Row 16
Set-TextCell 16 4 '80.76'
Set-TextCell 16 5 '  +2.22%  '

# Row 17
Set-TextCell 17 4 '1.013'
Set-TextCell 17 5 '  +0.75%  '

# Row 18
Set-TextCell 18 4 '0.000008844'
Set-TextCell 18 5 '  +2.27%  '

# Row 19
Set-TextCell 19 5 '  +0.46%  '

# Row 20
Set-TextCell 20 5 '  +2.19%  '

# Row 21
Set-TextCell 21 4 '26.948.80'
Set-TextCell 21 5 '  +3.60%  '

# Row 22
Set-TextCell 22 4 '5.194'
Set-TextCell 22 5 '  +4.41%  '

# Row 23
Set-TextCell 23 4 '10.92'
Set-TextCell 23 5 '  -1.05%  '

# Row 24
Set-TextCell 24 4 '2.286.49'
Set-TextCell 24 5 '  +17.07%  '

# Row 25
Set-TextCell 25 4 '154.00'
Set-TextCell 25 5 '  +1.41%  '

# Row 26
Set-TextCell 26 5 '  -1.10%  '

# Row 27
Set-TextCell 27 4 '18.37'
Set-TextCell 27 5 '  +1.48%  '

# Row 28
Set-TextCell 28 4 '5.237'
Set-TextCell 28 5 '  +4.29%  '

# Row 29
Set-TextCell 29 4 '1.913'
Set-TextCell 29 5 '  +14.27%  '

# Row 30
Set-TextCell 30 4 '114.71'
Set-TextCell 30 5 '  +0.18%  '

# Row 31
Set-TextCell 31 4 '0.08942'
Set-TextCell 31 5 '  +0.66%  '

# Row 32
Set-TextCell 32 4 '0.7414'
Set-TextCell 32 5 '  +3.08%  '

# Row 33
Set-TextCell 33 4 '1.157'
Set-TextCell 33 5 '  +5.16%  '

# Row 34
Set-TextCell 34 4 '4.435'
Set-TextCell 34 5 '  +3.20%  '

# Row 35
Set-TextCell 35 4 '2.805'
Set-TextCell 35 5 '  +2.40%  '

# Row 36
Set-TextCell 36 5 '  +0.66%  '

# Row 37
Set-TextCell 37 4 '1.115'
Set-TextCell 37 5 '  +4.69%  '

# Row 38
Set-TextCell 38 4 '0.05218'
Set-TextCell 38 5 '  +2.66%  '

# Row 39
Set-TextCell 39 5 '  +2.38%  '

# Row 40
Set-TextCell 40 5 '  +4.03%  '

# Row 41
Set-TextCell 41 4 '2.756'
Set-TextCell 41 5 '  +11.81%  '

# Row 42
Set-TextCell 42 4 '0.1648'
Set-TextCell 42 5 '  +3.21%  '

# Row 43
Set-TextCell 43 4 '6.463'
Set-TextCell 43 5 '  +5.34%  '

# Row 44
Set-TextCell 44 4 '8.279'
Set-TextCell 44 5 '  +4.05%  '

# Row 45
Set-TextCell 45 4 '107.29'
Set-TextCell 45 5 '  +2.61%  '

# Row 46
Set-TextCell 46 4 '10.35'
Set-TextCell 46 5 '  +3.99%  '

# Row 47
Set-TextCell 47 5 '  +0.54%  '

# Row 48
Set-TextCell 48 2 'Decentraland'
Set-TextCell 48 3 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextCell 48 4 '0.4573'
Set-TextCell 48 5 '  +2.64%  '

# Row 49
Set-TextCell 49 2 'NEARProtocol'
Set-TextCell 49 3 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 49 4 '1.645'
Set-TextCell 49 5 '  +5.14%  '

# Row 50
Set-TextCell 50 4 '0.06273'
Set-TextCell 50 5 '  +1.85%  '

# Row 51
Set-TextCell 51 4 '1.814'
Set-TextCell 51 5 '  +6.03%  '
